# chore: update ui import user
#
# The "comp_code" header was renamed to the fuller "company_code", and the
# newly-imported row gets blank placeholders in the (not-yet-populated)
# department/division/sub-division/level/position columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header A1: comp_code -> company_code
$ws.Range("A1").Value = "company_code"

# Row 2 (the freshly imported user) doesn't have org-structure data yet -
# leave those columns blank.
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""
